$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1727.9333
$ws.Range("I19").Value = 2004.2858
$ws.Range("J19").Value = 1486.125
$ws.Range("K19").Value = 2004.2858
$ws.Range("L19").Value = 1486.125
$ws.Range("M19").Value = -1829.2858
$ws.Range("N19").Value = -1836.125

$ws.Range("H80").Value = 58096.156
$ws.Range("I80").Value = 2800.7144
$ws.Range("K80").Value = 8402.143199999999
$ws.Range("M80").Value = -7404.143199999999

$ws.Range("H83").Value = 58096.156
$ws.Range("I83").Value = 2800.7144
$ws.Range("K83").Value = 25206.4296
$ws.Range("M83").Value = -20214.4296

$ws.Range("H112").Value = 1540.6538
$ws.Range("J112").Value = 1542.28
$ws.Range("L112").Value = 4626.84
$ws.Range("N112").Value = -6842.84

$ws.Range("H138").Value = 5368.0796
$ws.Range("I138").Value = 4638.357
$ws.Range("J138").Value = 5576.5713
$ws.Range("K138").Value = 13915.071
$ws.Range("L138").Value = 16729.7139
$ws.Range("M138").Value = -8775.071
$ws.Range("N138").Value = -27009.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5109.4
$ws.Range("I2").Value = 5275.5
$ws.Range("J2").Value = 4998.6665
$ws.Range("K2").Value = 5275.5
$ws.Range("L2").Value = 4998.6665
$ws.Range("M2").Value = -5162.5
$ws.Range("N2").Value = -5224.6665

$ws.Range("H45").Value = 1921.9
$ws.Range("I45").Value = 1865.1578
$ws.Range("K45").Value = 1865.1578
$ws.Range("M45").Value = -1488.1578

$ws.Range("H63").Value = 2999.5
$ws.Range("I63").Value = 2999.5
$ws.Range("K63").Value = 2999.5
$ws.Range("M63").Value = -2313.5

$ws.Range("H66").Value = 2999.5
$ws.Range("I66").Value = 2999.5
$ws.Range("K66").Value = 14997.5
$ws.Range("M66").Value = -11565.5

$ws.Range("H116").Value = 5109.4
$ws.Range("I116").Value = 5275.5
$ws.Range("J116").Value = 4998.6665
$ws.Range("K116").Value = 5275.5
$ws.Range("L116").Value = 4998.6665
$ws.Range("M116").Value = -2981.5
$ws.Range("N116").Value = -9586.666499999999

$ws.Range("H132").Value = 8999.571
$ws.Range("I132").Value = 3749.25
$ws.Range("K132").Value = 11247.75
$ws.Range("M132").Value = -8717.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5109.4
$ws.Range("I3").Value = 5275.5
$ws.Range("J3").Value = 4998.6665
$ws.Range("K3").Value = 5275.5
$ws.Range("L3").Value = 4998.6665
$ws.Range("M3").Value = -5161.5
$ws.Range("N3").Value = -5226.6665

$ws.Range("H86").Value = 3073
$ws.Range("I86").Value = 2335.25
$ws.Range("J86").Value = 7499.5
$ws.Range("K86").Value = 2335.25
$ws.Range("L86").Value = 7499.5
$ws.Range("M86").Value = -1212.25
$ws.Range("N86").Value = -9745.5

$ws.Range("H89").Value = 3073
$ws.Range("I89").Value = 2335.25
$ws.Range("J89").Value = 7499.5
$ws.Range("K89").Value = 11676.25
$ws.Range("L89").Value = 37497.5
$ws.Range("M89").Value = -6060.25
$ws.Range("N89").Value = -48729.5

$ws.Range("H92").Value = 59999.5
$ws.Range("J92").Value = 59999.5
$ws.Range("L92").Value = 59999.5
$ws.Range("N92").Value = -64991.5

$ws.Range("H105").Value = 3018.5
$ws.Range("I105").Value = 2100
$ws.Range("J105").Value = 3937
$ws.Range("K105").Value = 2100
$ws.Range("L105").Value = 3937
$ws.Range("M105").Value = -353
$ws.Range("N105").Value = -7431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12804
$ws.Range("I31").Value = 16919.5
$ws.Range("K31").Value = 16919.5
$ws.Range("M31").Value = -16624.5

$ws.Range("H34").Value = 12804
$ws.Range("I34").Value = 16919.5
$ws.Range("K34").Value = 16919.5
$ws.Range("M34").Value = -16717.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 852
$ws.Range("I122").Value = 704
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6336
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3886
$ws.Range("N122").Value = -13900

$ws.Range("H139").Value = 3359.077
$ws.Range("I139").Value = 3267.6
$ws.Range("J139").Value = 3664
$ws.Range("K139").Value = 9802.799999999999
$ws.Range("L139").Value = 10992
$ws.Range("M139").Value = -4662.799999999999
$ws.Range("N139").Value = -21272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11492.167
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 13240.6
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 13240.6
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -15236.6

$ws.Range("H83").Value = 11492.167
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 13240.6
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 66203
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -76187

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6629.4443
$ws.Range("I61").Value = 6510.154
$ws.Range("J61").Value = 6939.6
$ws.Range("K61").Value = 6510.154
$ws.Range("L61").Value = 6939.6
$ws.Range("M61").Value = -6308.154
$ws.Range("N61").Value = -7343.6

$ws.Range("H82").Value = 1876.6364
$ws.Range("I82").Value = 1545.3846
$ws.Range("J82").Value = 2355.111
$ws.Range("K82").Value = 1545.3846
$ws.Range("L82").Value = 2355.111
$ws.Range("M82").Value = -1184.3846
$ws.Range("N82").Value = -3077.111

$ws.Range("H85").Value = 1876.6364
$ws.Range("I85").Value = 1545.3846
$ws.Range("J85").Value = 2355.111
$ws.Range("K85").Value = 1545.3846
$ws.Range("L85").Value = 2355.111
$ws.Range("M85").Value = -297.3846000000001
$ws.Range("N85").Value = -4851.111

$ws.Range("H113").Value = 6629.4443
$ws.Range("I113").Value = 6510.154
$ws.Range("J113").Value = 6939.6
$ws.Range("K113").Value = 6510.154
$ws.Range("L113").Value = 6939.6
$ws.Range("M113").Value = -4340.154
$ws.Range("N113").Value = -11279.6

$ws.Range("H132").Value = 3768.4849
$ws.Range("I132").Value = 3443.1853
$ws.Range("K132").Value = 10329.5559
$ws.Range("M132").Value = -7799.555899999999

$ws.Range("H136").Value = 10335.667
$ws.Range("I136").Value = 13564.2
$ws.Range("K136").Value = 40692.60000000001
$ws.Range("M136").Value = -38142.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1250.25
$ws.Range("I81").Value = 1250.25
$ws.Range("K81").Value = 2500.5
$ws.Range("M81").Value = -1439.5

$ws.Range("H84").Value = 1250.25
$ws.Range("I84").Value = 1250.25
$ws.Range("K84").Value = 12502.5
$ws.Range("M84").Value = -7198.5

$ws.Range("H113").Value = 483.07693
$ws.Range("I113").Value = 360.375
$ws.Range("J113").Value = 679.4
$ws.Range("K113").Value = 1081.125
$ws.Range("L113").Value = 2038.2
$ws.Range("M113").Value = 1088.875
$ws.Range("N113").Value = -6378.2

$ws.Range("H119").Value = 70348.5
$ws.Range("J119").Value = 70348.5
$ws.Range("L119").Value = 70348.5
$ws.Range("N119").Value = -80024.5

$ws.Range("H132").Value = 5314.6206
$ws.Range("I132").Value = 5217.077
$ws.Range("K132").Value = 15651.231
$ws.Range("M132").Value = -13121.231
